# Applies the diff:
#  1. Merges the "c" + "arcaças de motores: 63 a 450" runs into one run
#     "carcaças de motores: 63 a 450" inside the last existing list item.
#  2. Appends 12 new "PargrafodaLista" (numId 4) list paragraphs at the
#     end of the document body (before the final sectPr).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: replace the visible (non-paragraph-mark) content of a Range
# with a raw WordprocessingML <w:p> fragment's runs, keeping the
# paragraph's own mark / pPr untouched.
# ---------------------------------------------------------------------
function Replace-ParagraphRuns($paragraph, $runsXml) {
    $pr = $paragraph.Range
    [int]$s = $pr.Start
    [int]$e = $pr.End
    # Exclude the trailing paragraph-mark character from the target range
    # so only the run content is replaced.
    $targetRange = $d.Range($s, $e - 1)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $targetRange.InsertXML($xml)
}

# ---------------------------------------------------------------------
# Helper: append a brand-new paragraph (given as a full <w:p>...</w:p>
# fragment) at the very end of the document body, after the current
# last paragraph and before sectPr.
# ---------------------------------------------------------------------
function Append-Paragraph($paragraphXml) {
    $last = $d.Paragraphs.Last
    $r = $last.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $newLast = $d.Paragraphs.Last
    $nr = $newLast.Range
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $paragraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $nr.InsertXML($xml)
    # InsertXML leaves a stray empty paragraph (carrying the old
    # paragraph mark) behind the inserted content; fold it away.
    $trailing = $d.Paragraphs.Last.Range
    [int]$s = $trailing.Start
    [int]$e = $trailing.End
    $full = $d.Range($s - 1, $e)
    $full.Delete()
}

# 1) Merge the "c" / "arcaças de motores: 63 a 450" runs -----------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Apenas para c*450*") {
        Replace-ParagraphRuns $p '<w:r><w:t xml:space="preserve">Apenas para </w:t></w:r><w:r><w:t>carcaças de motores: 63 a 450</w:t></w:r>'
        break
    }
}

# 2) Append the new list items --------------------------------------------
$pStyleOpen = '<w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>'

Append-Paragraph ('<w:p>' + $pStyleOpen + '<w:r><w:t>ótima resposta em relação</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>ao acelerômetro padrão dentro da faixa de frequências de trabalho de até 820 H</w:t></w:r></w:p>')

Append-Paragraph ('<w:p>' + $pStyleOpen + '<w:r><w:t>Interface para a visualização dos dados coletados</w:t></w:r></w:p>')

Append-Paragraph ('<w:p>' + $pStyleOpen + '<w:r><w:t>Interface para configuração do ambiente do usuário</w:t></w:r></w:p>')

Append-Paragraph ('<w:p>' + $pStyleOpen + '<w:r><w:t>Dispositivo discreto</w:t></w:r></w:p>')

Append-Paragraph ('<w:p>' + $pStyleOpen + '<w:r><w:t>Instalação intuitiva</w:t></w:r></w:p>')

Append-Paragraph ('<w:p>' + $pStyleOpen + '<w:r><w:t>Fixação via parafusos</w:t></w:r></w:p>')

Append-Paragraph ('<w:p>' + $pStyleOpen + '<w:r><w:t>Interface para gerenciamento de usuários</w:t></w:r></w:p>')

Append-Paragraph ('<w:p>' + $pStyleOpen + '<w:r><w:t>Pode ser instalado em motores acionados por inversor de frequência</w:t></w:r></w:p>')

Append-Paragraph ('<w:p>' + $pStyleOpen + '<w:r><w:t xml:space="preserve">Pode ser instalado mais de um Motor </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Scan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> em um mesmo motor</w:t></w:r></w:p>')

Append-Paragraph ('<w:p>' + $pStyleOpen + '<w:r><w:t>Medições agendadas ou em tempo real</w:t></w:r></w:p>')

Append-Paragraph ('<w:p>' + $pStyleOpen + '<w:r><w:t>Gerenciamento de sensores para medição</w:t></w:r></w:p>')

Append-Paragraph ('<w:p>' + $pStyleOpen + '<w:r><w:t xml:space="preserve">Imagem da planta de motores </w:t></w:r><w:r><w:t>editável</w:t></w:r></w:p>')

Write-Host ("Final paragraph count: " + $d.Paragraphs.Count)
